# "Loan RBI, Variable Instalments"
#
# On the "Repayment Schedule" sheet, a new (blank) column is inserted
# immediately before the old column N ("Late"). That pushes the old
# N/O/P columns ("Late", the blank spacer column, and "#"/Outstanding)
# one slot to the right, into O/P/Q, and widens the used range from
# A1:P14 to A1:Q14.
#
# The workbook's active sheet also moves from "Input" to
# "Repayment Schedule", with the selection there ending up on S4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column at N; existing N/O/P shift right to O/P/Q.
$ws.Columns("N").Insert()

# The newly inserted column N keeps a specific width (stored width 10,
# no "best fit") instead of the engine's default for a fresh column.
# ColumnWidth (character units) differs from the stored <col width>
# (character units + ~0.8333 padding), so back the padding out here.
$ws.Columns("N").ColumnWidth = 9.166666666666666

# Make "Repayment Schedule" the active sheet/tab, with S4 selected.
$ws.Activate()
$ws.Range("S4").Select()
